# Apply "Major bug fixes and some new features" update to the TableData sheet.
# Adds three new data-source rows (EW_Altersklassen, Familien, FlaecheDichte, Wanderungen... )
# Note: the original sheet already had two blank placeholder rows (row 3 & row 4) that get filled in,
# plus two brand-new rows (5 & 6) appended below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TableData")

# Style the two brand-new rows (5 & 6) the same as row 4 (center horizontal) before
# filling in values, so the fill order below drives the shared-string insertion order.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:D4").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122) | Out-Null

# Row 3: EW_Altersklassen (fill order matches original authoring: D, C, B, then A)
$ws.Cells.Item(3, 4).Value = "EW_Altersklassen_link"
$ws.Cells.Item(3, 3).Value = "EW_Altersklassen_link.xlsx"
$ws.Cells.Item(3, 2).Value = "EW Altersklassen"

# Row 4: Familien (fill order: D, C, B, A)
$ws.Cells.Item(4, 4).Value = "Familien_link"
$ws.Cells.Item(4, 3).Value = "Familien_link.xlsx"
$ws.Cells.Item(4, 2).Value = "Familien"
$ws.Cells.Item(4, 1).Value = "Familien"

$ws.Cells.Item(3, 1).Value = "EW_Altersklassen"

# Row 5: FlaecheDichte (new row)
$ws.Cells.Item(5, 1).Value = "FlaecheDichte"
$ws.Cells.Item(5, 2).Value = "Flächen und Dichten"
$ws.Cells.Item(5, 3).Value = "FlaecheDichte_link.xlsx"
$ws.Cells.Item(5, 4).Value = "FlaecheDichte_link"

# Row 6: Wanderungen (new row)
$ws.Cells.Item(6, 1).Value = "Wanderungen"
$ws.Cells.Item(6, 2).Value = "Wanderungen"
$ws.Cells.Item(6, 3).Value = "Wanderungen_link.xlsx"
$ws.Cells.Item(6, 4).Value = "Wanderungen_link"

$ws.Range("A1").Select() | Out-Null
$ws.Range("C3").Select() | Out-Null
